$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 121 (old rows 121:129 shift down to 123:131)
$ws.Rows("121:122").Insert()

# Row 121: Camote, 1a (guarda), new week data
$ws.Range("A121").Value = 11
$ws.Range("B121").Value = "Vega Monumental Concepción"
$ws.Range("C121").Value = "Bíobío"
$ws.Range("D121").Value = 44491
$ws.Range("E121").Value = 8
$ws.Range("F121").Value = 100112045
$ws.Range("G121").Value = "Zapallo"
$ws.Range("H121").Value = "Camote"
$ws.Range("I121").Value = "1a (guarda)"
$ws.Range("J121").Value = 600
$ws.Range("K121").Value = 550
$ws.Range("L121").Value = 600
$ws.Range("M121").Value = 575
$ws.Range("N121").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O121").Value = "Región de O'Higgins"
$ws.Range("P121").Value = 575
$ws.Range("Q121").Value = 1
$ws.Range("R121").Value = "Hortaliza"

# Row 122: Camote, 2a (guarda), new week data
$ws.Range("A122").Value = 11
$ws.Range("B122").Value = "Vega Monumental Concepción"
$ws.Range("C122").Value = "Bíobío"
$ws.Range("D122").Value = 44491
$ws.Range("E122").Value = 8
$ws.Range("F122").Value = 100112045
$ws.Range("G122").Value = "Zapallo"
$ws.Range("H122").Value = "Camote"
$ws.Range("I122").Value = "2a (guarda)"
$ws.Range("J122").Value = 300
$ws.Range("K122").Value = 450
$ws.Range("L122").Value = 450
$ws.Range("M122").Value = 450
$ws.Range("N122").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O122").Value = "Región de O'Higgins"
$ws.Range("P122").Value = 450
$ws.Range("Q122").Value = 1
$ws.Range("R122").Value = "Hortaliza"

# Ensure D121/D122 keep date style index 2 (same as other D cells) and numeric type
$ws.Range("D121:D122").NumberFormat = $ws.Range("D123").NumberFormat
